$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- C1 header: "habilidades_cnt_em" -> "habilidades", restyled like a title band ---
$c1 = $ws.Range("C1")
$c1.Value2 = "habilidades"
$c1.Font.Name = "Calibri"
$c1.Font.Size = 14
$c1.Font.Bold = $true
$c1.Font.Color = 16777215       # RGB(255,255,255) white
$c1.Interior.Color = 12611584   # RGB(0,112,192) blue
$c1.HorizontalAlignment = -4108 # xlCenter
$c1.VerticalAlignment = -4108   # xlCenter
$c1.WrapText = $true

# Row 1 shrinks now that the header text is shorter
$ws.Rows.Item(1).RowHeight = 37.5

# Print setup: portrait, paper size 9 (A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection moves to C1
$c1.Select() | Out-Null
